# Auto-generated edit script: update cryptos list (Price / Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.893.92'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '2.824.63'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'354.08"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'113.33"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.01%  '
$ws.Range('E7').Value = '  +4.65%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = "'0.604"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.75%  '
$ws.Range('D10').Value = "'41.62"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.37%  '
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').Value = "'19.92"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('D14').Value = "'7.76"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').Value = '3.267.46'
$ws.Range('E15').Value = '  +1.63%  '
$ws.Range('D16').Value = '2.827.49'
$ws.Range('E16').Value = '  +1.75%  '
$ws.Range('D17').Value = "'0.888"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').Value = '51.796.54'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = "'7.49"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.05%  '
$ws.Range('E20').Value = '  -2.79%  '
$ws.Range('D21').Value = "'13.44"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.90%  '
$ws.Range('D22').Value = '0.0₃0990'
$ws.Range('E22').Value = '  +1.65%  '
$ws.Range('D23').Value = "'270.16"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.67%  '
$ws.Range('D24').Value = "'69.77"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = "'2.79"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.18%  '
$ws.Range('D26').Value = "'26.78"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('E30').Value = '  -1.88%  '
$ws.Range('D31').Value = "'50.77"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('D32').Value = "'34.00"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.61%  '
$ws.Range('D33').Value = "'0.0454"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +26.19%  '
$ws.Range('D34').Value = "'5.86"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.88%  '
$ws.Range('D35').Value = "'5.27"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.16%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').Value = "'0.999"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').Value = "'2.07"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('E39').Value = '  -1.52%  '
$ws.Range('D40').Value = "'18.24"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.62%  '
$ws.Range('D41').Value = "'23.86"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.36%  '
$ws.Range('E42').Value = '  +2.24%  '
$ws.Range('D43').Value = "'126.24"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('D44').Value = "'2.52"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.85%  '
$ws.Range('D45').Value = "'2.30"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.97%  '
$ws.Range('D46').Value = '2.092.86'
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('D47').Value = "'3.36"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.88%  '
$ws.Range('D48').Value = "'2.29"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.86%  '
$ws.Range('D49').Value = "'5.68"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.60%  '
$ws.Range('D50').Value = "'0.936"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.11%  '
$ws.Range('D51').Value = "'60.88"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.19%  '
